# Remove the "Employee ID:" label and its placeholder value ("#########")
# from the employee-info table, leaving the two table cells empty.
$d = $word.ActiveDocument

$d.Content.Find.Execute("Employee ID:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

$d.Content.Find.Execute("#########", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)
